$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 33: re-enter the header-style legend values (2,3,4,5) under S:V,
#     mirroring row 3's S3:V3 values (same unstyled cells). ---
$ws.Range("S33").Value = 2
$ws.Range("T33").Value = 3
$ws.Range("U33").Value = 4
$ws.Range("V33").Value = 5

# --- Row 35: totals row summing each of the S:V "match" columns over the
#     student rows (4:31). S35 gets its own formula; T35:V35 are entered as
#     one multi-cell formula so Excel records them as a shared formula group
#     (same pattern already used throughout rows 4:31 for these columns). ---
$ws.Range("S35").Formula = "=SUM(S4:S31)"
$ws.Range("T35:V35").Formula = "=SUM(T4:T31)"

# --- View state: move the active selection of the bottom-right (scrollable)
#     pane down to row 28, and bring the frozen pane's visible window down
#     so row 11 is at the top (same freeze boundary at B/row3 is preserved). ---
$win = $excel.ActiveWindow
$win.ScrollRow = 11
$win.ScrollColumn = 3
$ws.Range("X28").Select() | Out-Null
